$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Remove the Topic ("Array") value from B24 - fully clear cell (content + formatting)
# so the cell falls back to the row's default style and is omitted from the XML,
# matching the rest of the sparse rows below it.
$ws.Range("B24").Clear()

# Fill in the serial numbers for the newly-populated rows (25-33)
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32

# Populate row 33 with the new LeetCode problem entry
$ws.Range("F33").Value = "https://leetcode.com/problems/intersection-of-two-arrays/description/"
$ws.Range("G33").Value = 45557
$ws.Range("G33").NumberFormat = "d-mmm"
$ws.Range("H33").Value = "Sept"

# Match the final on-screen selection from the edit session
$ws.Range("A33").Select() | Out-Null

